$wb = $excel.ActiveWorkbook

# Rename the second sheet from "TestCase01" to "signUpWithPhoneNumber"
$ws = $wb.Worksheets.Item(2)
$ws.Name = "signUpWithPhoneNumber"

# Populate the sheet with header and data rows
$ws.Range("A1").Value = "COUNTRY_CODE"
$ws.Range("B1").Value = "OTP"

$ws.Range("A2").Value = "IN"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "1111"
$ws.Range("B2").Style = "Normal"

$ws.Range("A3").Value = "US"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "1111"
$ws.Range("B3").Style = "Normal"
